# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values for the two handback rows
# on the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-27 02:37:08"
$wsZhCn.Range("G2").Value = "2016-01-27 02:37:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-27 02:37:19"
$wsDeDe.Range("G2").Value = "2016-01-27 02:38:19"
